$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1928733333333333
$ws.Range("H2").Value = 0.57862
$ws.Range("I2").Value = 0.2550396805282215
$ws.Range("J2").Value = 0.2550396805282216
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 0.1434900450666667
$ws.Range("R2").Value = 1.2914104056
$ws.Range("S2").Value = 0.001596294180838868
$ws.Range("T2").Value = 0.001596294180838868

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1928733333333333
$ws.Range("H3").Value = 0.57862
$ws.Range("I3").Value = 0.2550396805282215
$ws.Range("J3").Value = 0.2550396805282216
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 17.00174292985778
$ws.Range("R3").Value = 153.01568636872
$ws.Range("S3").Value = 0.189140530901924
$ws.Range("T3").Value = 0.189140530901924

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1928733333333333
$ws.Range("H4").Value = 0.57862
$ws.Range("I4").Value = 0.2550396805282215
$ws.Range("J4").Value = 0.2550396805282216
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 5.741568982084444
$ws.Range("R4").Value = 51.67412083876
$ws.Range("S4").Value = 0.06387365165805124
$ws.Range("T4").Value = 0.06387365165805126

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1928733333333333
$ws.Range("H5").Value = 0.57862
$ws.Range("I5").Value = 0.2550396805282215
$ws.Range("J5").Value = 0.2550396805282216
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.03858090290444445
$ws.Range("R5").Value = 0.3472281261400001
$ws.Range("S5").Value = 0.0004292037874074154
$ws.Range("T5").Value = 0.0004292037874074155

$ws.Range("G6").Value = 0.5633750000000001
$ws.Range("H6").Value = 1.690125
$ws.Range("I6").Value = 0.7449603194717784
$ws.Range("J6").Value = 0.7449603194717785
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 0.419128465
$ws.Range("R6").Value = 3.772156185
$ws.Range("S6").Value = 0.004662709035965386
$ws.Range("T6").Value = 0.004662709035965387

$ws.Range("G7").Value = 0.5633750000000001
$ws.Range("H7").Value = 1.690125
$ws.Range("I7").Value = 0.7449603194717784
$ws.Range("J7").Value = 0.7449603194717785
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("Q7").Value = 49.66138531216667
$ws.Range("R7").Value = 446.9524678095
$ws.Range("S7").Value = 0.5524716390560546
$ws.Range("T7").Value = 0.5524716390560547

$ws.Range("G8").Value = 0.5633750000000001
$ws.Range("H8").Value = 1.690125
$ws.Range("I8").Value = 0.7449603194717784
$ws.Range("J8").Value = 0.7449603194717785
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 16.77088464941667
$ws.Range("R8").Value = 150.93796184475
$ws.Range("S8").Value = 0.1865722849340912
$ws.Range("T8").Value = 0.1865722849340913

$ws.Range("G9").Value = 0.5633750000000001
$ws.Range("H9").Value = 1.690125
$ws.Range("I9").Value = 0.7449603194717784
$ws.Range("J9").Value = 0.7449603194717785
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 0.1126932157916667
$ws.Range("R9").Value = 1.014238942125
$ws.Range("S9").Value = 0.001253686445667205
$ws.Range("T9").Value = 0.001253686445667205
